$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# "Generate Report for Handoff" - refresh the localization-status report for
# the b47a380c-7b15-4457-ab9f-3b63fe9788b4.md file: a new handoff round was
# generated (zh-cn is ready for handoff; de-de's existing handback turned out
# to be stale against the latest source, so an error note is attached).
# ---------------------------------------------------------------------------

$errorDetail = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/bf518601f9dae2ca4d8e9033cff490895c4a83fe/e2e/b47a380c-7b15-4457-ab9f-3b63fe9788b4.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/5b44ed7632355ca31c7832aba46367a12de705d4/e2e/b47a380c-7b15-4457-ab9f-3b63fe9788b4.md."

# --- Overview sheet ---------------------------------------------------------
$ov = $wb.Worksheets.Item("Overview")
$ov.Range("E3").Value = "Ready for handoff"
$ov.Range("F3").Value = "Ready for handoff"
$ov.Range("G3").Value = "2016-09-07 05:02:51"

# --- zh-cn sheet -------------------------------------------------------------
$zh = $wb.Worksheets.Item("zh-cn")
$zh.Range("C3").Value = "Ready for handoff"
$zh.Range("H3").Value = "2016-09-07 05:02:45"
$zh.Range("P3").Value = $errorDetail
$zh.Range("P1").ColumnWidth = 39.15

# --- de-de sheet -------------------------------------------------------------
$de = $wb.Worksheets.Item("de-de")
$de.Range("H3").Value = "2016-09-07 05:02:51"
$de.Range("P3").Value = $errorDetail
$de.Range("P1").ColumnWidth = 39.15
